# ==== Auto-generated Excel COM-interop edit script ====
$wb = $excel.ActiveWorkbook

# --- 1) Insert new worksheet 'Đơn phụ phẫu 1' between 'Đơn sale chính' and 'Lương' ---
$sheetChinh = $wb.Worksheets.Item(1)
$newSheet = $wb.Worksheets.Add($null, $sheetChinh)
$newSheet.Name = "Đơn phụ phẫu 1"

# Header row 1 (same columns/labels as 'Đơn sale chính')
$newSheet.Range("A1").Value = "Tiền tố"
$newSheet.Range("B1").Value = "Mã dịch vụ"
$newSheet.Range("C1").Value = "Ngày thực hiện"
$newSheet.Range("D1").Value = "Cơ sở"
$newSheet.Range("E1").Value = "Khách hàng"
$newSheet.Range("F1").Value = "Nguồn khách"
$newSheet.Range("G1").Value = "Nhóm dịch vụ"
$newSheet.Range("H1").Value = "Tên dịch vụ"
$newSheet.Range("I1").Value = "Sale chính"
$newSheet.Range("J1").Value = "Đơn giá gốc"
$newSheet.Range("K1").Value = "Sale phụ"
$newSheet.Range("L1").Value = "Upsale"
$newSheet.Range("M1").Value = "Đơn giá"
$newSheet.Range("N1").Value = "Thanh toán lần đầu"
$newSheet.Range("O1").Value = "Trả sau"
$newSheet.Range("P1").Value = "Đã thanh toán"
$newSheet.Range("Q1").Value = "Dư nợ"
$newSheet.Range("R1").Value = "Bác sĩ 1"
$newSheet.Range("S1").Value = "Bác sĩ 2"
$newSheet.Range("T1").Value = "Phụ phẫu 1"
$newSheet.Range("U1").Value = "Phụ phẫu 2"
$newSheet.Range("V1").Value = "Công phụ phẫu 1"
$newSheet.Range("W1").Value = "Công phụ phẫu 2"
$newSheet.Range("X1").Value = "Tỉ lệ chiết khấu sale chính"
$newSheet.Range("Y1").Value = "Tỉ lệ chiết khấu sale phụ"
$newSheet.Range("Z1").Value = "Chiết khấu sale chính"
$newSheet.Range("AA1").Value = "Chiết khấu sale phụ"

# Row 2 - data
$newSheet.Range("A2").Value = "HD-LUXURY"
$newSheet.Range("B2").Value = 555
$newSheet.Range("C2").NumberFormat = "@"
$newSheet.Range("C2").Value = "07-16-2024"
$newSheet.Range("D2").Value = "LONG XUYÊN"
$newSheet.Range("E2").Value = "Nguyễn thị hồng cúc"
$newSheet.Range("F2").Value = "Cá nhân"
$newSheet.Range("G2").Value = "Tiêm"
$newSheet.Range("H2").Value = "Tiêm Filler"
$newSheet.Range("I2").Value = "Lê Văn Linh"
$newSheet.Range("J2").Value = 7000000
$newSheet.Range("M2").Value = 7000000
$newSheet.Range("N2").Value = 7000000
$newSheet.Range("O2").Value = 0
$newSheet.Range("P2").Value = 7000000
$newSheet.Range("Q2").Value = 0
$newSheet.Range("R2").Value = "Đặng Ngọc Mai"
$newSheet.Range("T2").Value = "Sang sang"
$newSheet.Range("V2").Value = 50000
$newSheet.Range("W2").Value = 0
$newSheet.Range("X2").Value = 0.13
$newSheet.Range("Y2").Value = 0
$newSheet.Range("Z2").Value = 910000
$newSheet.Range("AA2").Value = 0

# Row 3 - totals
$newSheet.Range("A3").Value = "Tổng"
$newSheet.Range("B3").Value = 1
$newSheet.Range("J3").Value = 7000000
$newSheet.Range("L3").Value = 0
$newSheet.Range("M3").Value = 7000000
$newSheet.Range("N3").Value = 7000000
$newSheet.Range("O3").Value = 0
$newSheet.Range("P3").Value = 7000000
$newSheet.Range("Q3").Value = 0
$newSheet.Range("V3").Value = 50000
$newSheet.Range("W3").Value = 0
$newSheet.Range("X3").Value = 0.13
$newSheet.Range("Y3").Value = 0
$newSheet.Range("Z3").Value = 910000
$newSheet.Range("AA3").Value = 0

# --- 2) 'Đơn sale chính' row2 G2: 'Môi' -> 'Tiêm' ---
$sheetChinh.Range("G2").Value = "Tiêm"

# --- 3) Update 'Lương' sheet: shift/insert rows and new totals ---
$luong = $wb.Worksheets.Item("Lương")
$luong.Range("A1").Value = "Danh mục"
$luong.Range("B1").Value = 6
$luong.Range("A2").Value = "Ngày công"
$luong.Range("B2").Value = 17
$luong.Range("A3").Value = "Phụ cấp"
$luong.Range("B3").Value = 595000
$luong.Range("A4").Value = "Lương cơ bản tại CẦN THƠ"
$luong.Range("B4").Value = ""
$luong.Range("A5").Value = "Chiết khấu sale chính tại CẦN THƠ"
$luong.Range("B5").Value = 0
$luong.Range("A6").Value = "Chiết khấu sale phụ tại CẦN THƠ"
$luong.Range("B6").Value = 0
$luong.Range("A7").Value = "Đơn 1 bác sĩ tại CẦN THƠ"
$luong.Range("B7").Value = 0
$luong.Range("A8").Value = "Đơn 2 bác sĩ tại CẦN THƠ"
$luong.Range("B8").Value = 0
$luong.Range("A9").Value = "Công phụ phẫu 1 tại CẦN THƠ"
$luong.Range("B9").Value = 0
$luong.Range("A10").Value = "Công phụ phẫu 2 tại CẦN THƠ"
$luong.Range("B10").Value = 0
$luong.Range("A11").Value = "Ứng lương tại CẦN THƠ"
$luong.Range("B11").Value = 0
$luong.Range("A12").Value = "Lương cơ bản tại LONG XUYÊN"
$luong.Range("B12").Value = 1821428.571428571
$luong.Range("A13").Value = "Chiết khấu sale chính tại LONG XUYÊN"
$luong.Range("B13").Value = 150000
$luong.Range("A14").Value = "Chiết khấu sale phụ tại LONG XUYÊN"
$luong.Range("B14").Value = 0
$luong.Range("A15").Value = "Đơn 1 bác sĩ tại LONG XUYÊN"
$luong.Range("B15").Value = 0
$luong.Range("A16").Value = "Đơn 2 bác sĩ tại LONG XUYÊN"
$luong.Range("B16").Value = 0
$luong.Range("A17").Value = "Công phụ phẫu 1 tại LONG XUYÊN"
$luong.Range("B17").Value = 50000
$luong.Range("A18").Value = "Công phụ phẫu 2 tại LONG XUYÊN"
$luong.Range("B18").Value = 0
$luong.Range("A19").Value = "Ứng lương tại LONG XUYÊN"
$luong.Range("B19").Value = 0
$luong.Range("A20").Value = "Lương cơ bản tại SÓC TRĂNG"
$luong.Range("B20").Value = ""
$luong.Range("A21").Value = "Chiết khấu sale chính tại SÓC TRĂNG"
$luong.Range("B21").Value = 0
$luong.Range("A22").Value = "Chiết khấu sale phụ tại SÓC TRĂNG"
$luong.Range("B22").Value = 0
$luong.Range("A23").Value = "Đơn 1 bác sĩ tại SÓC TRĂNG"
$luong.Range("B23").Value = 0
$luong.Range("A24").Value = "Đơn 2 bác sĩ tại SÓC TRĂNG"
$luong.Range("B24").Value = 0
$luong.Range("A25").Value = "Công phụ phẫu 1 tại SÓC TRĂNG"
$luong.Range("B25").Value = 0
$luong.Range("A26").Value = "Công phụ phẫu 2 tại SÓC TRĂNG"
$luong.Range("B26").Value = 0
$luong.Range("A27").Value = "Ứng lương tại SÓC TRĂNG"
$luong.Range("B27").Value = 0
$luong.Range("A28").Value = "Tổng lương tại CẦN THƠ"
$luong.Range("B28").Value = 0
$luong.Range("A29").Value = "Tổng lương tại LONG XUYÊN"
$luong.Range("B29").Value = 7078928.571428571
$luong.Range("A30").Value = "Tổng lương tại SÓC TRĂNG"
$luong.Range("B30").Value = 0
$luong.Range("A31").Value = "Tổng lương"
$luong.Range("B31").Value = 7078928.571428571

